# Update countries & provincias Spain
# Applies the COVID dataset refresh: several countries' case numbers changed,
# and three country pairs swapped places in the (descending, by total cases)
# sorted list because one country's total overtook its neighbour's. The
# timestamp banner in A1 is bumped to the new refresh time.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rusia (row 7) : values refreshed, country stays in place ---
$ws.Range("B7").Value = 845443
$ws.Range("C7").Value = 5462
$ws.Range("D7").Value = 646524
$ws.Range("E7").Value = 184861
$ws.Range("F7").Value = 0
$ws.Range("G7").Value = 95
$ws.Range("H7").Value = 14058

# --- Israel / Ucrania (rows 36-37) swap order: Ucrania's total (71056)
#     overtakes Israel's (70970), so Ucrania moves up to row 36 with its
#     refreshed figures and Israel drops to row 37 unchanged. ---
$ws.Range("A36").Value = "Ucrania"
$ws.Range("B36").Value = 71056
$ws.Range("C36").Value = 1172
$ws.Range("D36").Value = 39308
$ws.Range("E36").Value = 30039
$ws.Range("F36").Value = 0
$ws.Range("G36").Value = 16
$ws.Range("H36").Value = 1709

$ws.Range("A37").Value = "Israel"
$ws.Range("B37").Value = 70970
$ws.Range("C37").Value = 0
$ws.Range("D37").Value = 43850
$ws.Range("E37").Value = 26608
$ws.Range("F37").Value = 0
$ws.Range("G37").Value = 0
$ws.Range("H37").Value = 512

# --- Singapur (row 45) : values refreshed, country stays in place ---
$ws.Range("B45").Value = 52512
$ws.Range("C45").Value = 307
$ws.Range("D45").Value = 46491
$ws.Range("E45").Value = 5994
$ws.Range("F45").Value = 0
$ws.Range("G45").Value = 0

# --- Sudan del Sur / Eslovaquia (rows 123-124) swap order: Eslovaquia's
#     total (2337) overtakes Sudan del Sur's (2322), so Eslovaquia moves up
#     to row 123 with its refreshed figures and Sudan del Sur drops to row
#     124 unchanged. ---
$ws.Range("A123").Value = "Eslovaquia"
$ws.Range("B123").Value = 2337
$ws.Range("C123").Value = 45
$ws.Range("D123").Value = 1742
$ws.Range("E123").Value = 566
$ws.Range("F123").Value = 0
$ws.Range("G123").Value = 0
$ws.Range("H123").Value = 29

$ws.Range("A124").Value = "Sudan del Sur"
$ws.Range("B124").Value = 2322
$ws.Range("C124").Value = 0
$ws.Range("D124").Value = 1175
$ws.Range("E124").Value = 1101
$ws.Range("F124").Value = 0
$ws.Range("G124").Value = 0
$ws.Range("H124").Value = 46

# --- Lituania (row 127) : values refreshed, country stays in place ---
$ws.Range("B127").Value = 2093
$ws.Range("C127").Value = 18
$ws.Range("D127").Value = 1644
$ws.Range("E127").Value = 369
$ws.Range("F127").Value = 0
$ws.Range("G127").Value = 0

# --- Estonia (row 128) : values refreshed, country stays in place ---
$ws.Range("B128").Value = 2072
$ws.Range("C128").Value = 8
$ws.Range("D128").Value = 1934
$ws.Range("E128").Value = 69
$ws.Range("F128").Value = 0
$ws.Range("G128").Value = 0

# --- Barbados / Islas Turcas y Caicos (rows 187-188) swap order: Islas
#     Turcas y Caicos's total (114) overtakes Barbados's (110), so it moves
#     up to row 187 with its refreshed figures and Barbados drops to row 188
#     unchanged. ---
$ws.Range("A187").Value = "Islas Turcas y Caicos"
$ws.Range("B187").Value = 114
$ws.Range("C187").Value = 7
$ws.Range("D187").Value = 38
$ws.Range("E187").Value = 74
$ws.Range("F187").Value = 0
$ws.Range("G187").Value = 0
$ws.Range("H187").Value = 2

$ws.Range("A188").Value = "Barbados"
$ws.Range("B188").Value = 110
$ws.Range("C188").Value = 0
$ws.Range("D188").Value = 96
$ws.Range("E188").Value = 7
$ws.Range("F188").Value = 0
$ws.Range("G188").Value = 0
$ws.Range("H188").Value = 7

# --- Refresh banner timestamp ---
$ws.Range("A1").Value = "Datos actualizados a 1 de Agosto de 2020 a las 10:24"
